# "Generate Report for Handback"
#
# The handback pipeline ran for both locales (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Each row now also reports a "Latest Target File" (the source .md that was
#     handed back) and a "Latest Handback File" (the translated .xlf), mirroring
#     the existing "Latest Handoff File" columns (F = Latest Target File,
#     G = Latest Handback File in the per-locale tables).
#   - "Latest Handback DateTime" moves off the 0001-01-01 sentinel onto the
#     actual handback timestamp (per locale).

$wb = $excel.ActiveWorkbook

$sourceRepo = "https://github.com/OpenLocalizationTest/oltest/blob/69550adeff2b454a48ac6ae691dab72d724cae40/e2e"

$locales = @(
    @{
        Sheet = "zh-cn"
        HandoffRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e543ea9549ed2b05f48432dd2127c9c11a8df43f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ht"
        HandbackDateTime = "2016-03-09 11:34:04"
        Rows = @(
            @{ Row = 2; Base = "777a8731-c636-414f-9322-10137bdedc81"; Xlf = "777a8731-c636-414f-9322-10137bdedc81.ac3da6758951b4125e261d5d2cd69da5ed167cee.zh-cn.xlf" },
            @{ Row = 3; Base = "ec12cc59-17a4-4245-bfb3-c754a2cf89bf"; Xlf = "ec12cc59-17a4-4245-bfb3-c754a2cf89bf.b3d06500a6268c0c35d3b8f530b1fe784f3a4446.zh-cn.xlf" }
        )
    },
    @{
        Sheet = "de-de"
        HandoffRepo = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4aa56d31aa17fb0b0f4d4ea68034226425bae7aa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ht"
        HandbackDateTime = "2016-03-09 11:34:20"
        Rows = @(
            @{ Row = 2; Base = "777a8731-c636-414f-9322-10137bdedc81"; Xlf = "777a8731-c636-414f-9322-10137bdedc81.ac3da6758951b4125e261d5d2cd69da5ed167cee.de-de.xlf" },
            @{ Row = 3; Base = "ec12cc59-17a4-4245-bfb3-c754a2cf89bf"; Xlf = "ec12cc59-17a4-4245-bfb3-c754a2cf89bf.b3d06500a6268c0c35d3b8f530b1fe784f3a4446.de-de.xlf" }
        )
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    foreach ($row in $locale.Rows) {
        $r = $row.Row
        $mdName = $row.Base + ".md"
        $mdUrl = $sourceRepo + "/" + $mdName
        $xlfName = $row.Xlf
        $xlfUrl = $locale.HandoffRepo + "/" + $xlfName

        # Status: handed back, in sync with en-US
        $ws.Range("C" + $r).Value = "Handed back: in sync with en-US"

        # F = Latest Target File (the source markdown that was handed back)
        $ws.Hyperlinks.Add($ws.Range("F" + $r), $mdUrl, "", "", $mdName)

        # G = Latest Handback File (the translated .xlf)
        $ws.Hyperlinks.Add($ws.Range("G" + $r), $xlfUrl, "", "", $xlfName)

        # H = Latest Handback DateTime
        $ws.Range("H" + $r).Value = $locale.HandbackDateTime
    }
}
